$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 139
$ws1.Range("F5").Value = 2961
$ws1.Range("F6").Value = 299
$ws1.Range("F7").Value = 401

# Sheet "全部类型" (sheet4): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 139
$ws4.Range("F5").Value = 2961
$ws4.Range("F6").Value = 299
$ws4.Range("F9").Value = 401
